$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2962962962962963
$ws.Range("C2").Value = 0.4074074074074074
$ws.Range("P2").Value = 0.2592592592592592
$ws.Range("S2").Value = 0.03703703703703703
$ws.Range("C3").Value = 0.09090909090909091
$ws.Range("P3").Value = 0.7272727272727273
$ws.Range("S3").Value = 0.1818181818181818
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.5
$ws.Range("B6").Value = 0.08333333333333333
$ws.Range("J6").Value = 0.3333333333333333
$ws.Range("Q6").Value = 0.25
$ws.Range("R6").Value = 0.08333333333333333
$ws.Range("S6").Value = 0.25
$ws.Range("B7").Value = 0.1818181818181818
$ws.Range("J7").Value = 0.1818181818181818
$ws.Range("R7").Value = 0.2727272727272727
$ws.Range("S7").Value = 0.3636363636363636
$ws.Range("B8").Value = 0.1891891891891892
$ws.Range("D8").Value = 0.02702702702702703
$ws.Range("F8").Value = 0.02702702702702703
$ws.Range("J8").Value = 0.2702702702702703
$ws.Range("R8").Value = 0.05405405405405406
$ws.Range("S8").Value = 0.4324324324324325
$ws.Range("J9").Value = 0.125
$ws.Range("Q9").Value = 0.375
$ws.Range("R9").Value = 0.125
$ws.Range("S9").Value = 0.375
$ws.Range("B10").Value = 0.07627118644067797
$ws.Range("D10").Value = 0.0423728813559322
$ws.Range("F10").Value = 0.0423728813559322
$ws.Range("J10").Value = 0.09322033898305085
$ws.Range("O10").Value = 0.01694915254237288
$ws.Range("Q10").Value = 0.2627118644067797
$ws.Range("R10").Value = 0.07627118644067797
$ws.Range("S10").Value = 0.3898305084745763
$ws.Range("G11").Value = 0.125
$ws.Range("K11").Value = 0.1875
$ws.Range("L11").Value = 0.625
$ws.Range("S11").Value = 0.0625
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5
$ws.Range("H15").Value = 0.05882352941176471
$ws.Range("I15").Value = 0.05882352941176471
$ws.Range("J15").Value = 0.6470588235294118
$ws.Range("M15").Value = 0.05882352941176471
$ws.Range("S15").Value = 0.1764705882352941
$ws.Range("H16").Value = 0.3888888888888889
$ws.Range("J16").Value = 0.4444444444444444
$ws.Range("K16").Value = 0.05555555555555555
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.05555555555555555
$ws.Range("F17").Value = 0.05405405405405406
$ws.Range("H17").Value = 0.1351351351351351
$ws.Range("I17").Value = 0.02702702702702703
$ws.Range("J17").Value = 0.5945945945945946
$ws.Range("K17").Value = 0.02702702702702703
$ws.Range("O17").Value = 0.02702702702702703
$ws.Range("S17").Value = 0.1351351351351351
$ws.Range("H18").Value = 0.1875
$ws.Range("I18").Value = 0.125
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.0625
$ws.Range("O18").Value = 0.0625
$ws.Range("S18").Value = 0.0625
$ws.Range("F19").Value = 0.0202020202020202
$ws.Range("H19").Value = 0.2121212121212121
$ws.Range("I19").Value = 0.04040404040404041
$ws.Range("J19").Value = 0.404040404040404
$ws.Range("K19").Value = 0.1111111111111111
$ws.Range("M19").Value = 0.0101010101010101
$ws.Range("O19").Value = 0.09090909090909091
$ws.Range("S19").Value = 0.1111111111111111
